$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-108 is updated from 2023-10-05 (45204)
# to 2023-10-08 (45207) for every data row in the sheet.
for ($r = 2; $r -le 108; $r++) {
    $ws.Cells.Item($r, 3).Value = 45207
}
